$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the "Price" cells that would otherwise be auto-converted to numbers
# by Excel (single-dot decimal-looking strings) as Text, to preserve them
# as literal strings exactly like the source data (e.g. "1.000", "0.4740").
$ws.Range("D4:D12").NumberFormat = "@"
$ws.Range("D14:D16").NumberFormat = "@"
$ws.Range("D18:D20").NumberFormat = "@"
$ws.Range("D22:D28").NumberFormat = "@"
$ws.Range("D31:D36").NumberFormat = "@"
$ws.Range("D38:D43").NumberFormat = "@"
$ws.Range("D45:D51").NumberFormat = "@"

$ws.Range("B12").Value = "TRON"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("B19").Value = "Dai"
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D2").Value = "30.578.14"
$ws.Range("D3").Value = "1.923.77"
$ws.Range("D4").Value = "1.000"
$ws.Range("D5").Value = "247.31"
$ws.Range("D6").Value = "1.000"
$ws.Range("D7").Value = "0.4740"
$ws.Range("D8").Value = "0.2911"
$ws.Range("D9").Value = "0.06795"
$ws.Range("D10").Value = "105.41"
$ws.Range("D11").Value = "18.45"
$ws.Range("D12").Value = "0.07724"
$ws.Range("D13").Value = "1.907.69"
$ws.Range("D14").Value = "5.320"
$ws.Range("D15").Value = "0.6718"
$ws.Range("D16").Value = "287.64"
$ws.Range("D17").Value = "30.614.70"
$ws.Range("D18").Value = "0.000007631"
$ws.Range("D19").Value = "1.001"
$ws.Range("D20").Value = "12.96"
$ws.Range("D21").Value = "2.160.55"
$ws.Range("D22").Value = "5.450"
$ws.Range("D23").Value = "1.0000"
$ws.Range("D24").Value = "6.316"
$ws.Range("D25").Value = "9.412"
$ws.Range("D26").Value = "168.13"
$ws.Range("D27").Value = "20.86"
$ws.Range("D28").Value = "2.148"
$ws.Range("D31").Value = "4.201"
$ws.Range("D32").Value = "4.123"
$ws.Range("D33").Value = "0.05048"
$ws.Range("D34").Value = "0.7425"
$ws.Range("D35").Value = "1.160"
$ws.Range("D36").Value = "0.02077"
$ws.Range("D38").Value = "2.691"
$ws.Range("D39").Value = "2.068"
$ws.Range("D40").Value = "111.43"
$ws.Range("D41").Value = "0.8809"
$ws.Range("D42").Value = "5.960"
$ws.Range("D43").Value = "0.4372"
$ws.Range("D45").Value = "67.37"
$ws.Range("D46").Value = "7.278"
$ws.Range("D47").Value = "9.322"
$ws.Range("D48").Value = "48.08"
$ws.Range("D49").Value = "0.1237"
$ws.Range("D50").Value = "35.13"
$ws.Range("D51").Value = "0.4060"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +5.08%  "
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +3.41%  "
$ws.Range("E9").Value = "  +5.88%  "
$ws.Range("E10").Value = "  +9.01%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("E14").Value = "  +6.78%  "
$ws.Range("E15").Value = "  +5.58%  "
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("E18").Value = "  +3.62%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("E22").Value = "  +9.02%  "
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("E24").Value = "  +4.32%  "
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("E26").Value = "  +2.00%  "
$ws.Range("E27").Value = "  +8.07%  "
$ws.Range("E28").Value = "  +11.71%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +2.80%  "
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("E32").Value = "  +8.39%  "
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +2.49%  "
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("E36").Value = "  +8.15%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("E42").Value = "  +6.05%  "
$ws.Range("E43").Value = "  +8.10%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  +3.12%  "
$ws.Range("E46").Value = "  +3.22%  "
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("E48").Value = "  +17.38%  "
$ws.Range("E49").Value = "  +4.28%  "
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("E51").Value = "  +9.03%  "
